$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Weihnachtsmarkt_Bilder")

# Fill in newly-researched image filenames / attributions for four rows:
#  id 3   - "17. Berliner Weihnachtszeit vor dem Roten Rathaus"
#  id 15  - "Weihnachtsmarkt Hallen am Borsigturm"
#  id 174 - "Weihnachtsmarkt an der Kaiser-Wilhelm-Gedaechtniskirche" (Breitscheidplatz)
#  id 153 - "Umwelt- und Weihnachtsmarkt in der Sophienstrasse"
$ws.Range("C45").Value = "Christmas_market_Rotes_Rathaus_Berlin.jpg"
$ws.Range("D45").Value = "Leonhard Lenz, CC Zero, Public Domain Dedication"

$ws.Range("C48").Value = "Eingang_HallenAmBorsigturm.jpg"
$ws.Range("D48").Value = "Aiken Hartenfels, CC Attribution-Share Alike 3.0 de"

$ws.Range("C55").Value = "weihnachtsmarkt_Breitscheidplatz.jpg"
$ws.Range("D55").Value = "Ralf Roletschek, GFDL 1.2 via Wikimedia Commons"

$ws.Range("C60").Value = "Christmas_market_Sophienstraße_Berlin.jpg"
$ws.Range("D60").Value = "Leonhard Lenz, CC Zero, Public Domain Dedication"

# Strip the (unused) extra number-format style from the whole data range,
# leaving every data cell using the default/unstyled cell format
$ws.Range("B2:D64").ClearFormats()

# Update the sheet's active selection to match the author's last on-screen state
$ws.Range("F58").Select()
